# "Add files via upload" — vaccinesfull.xlsx
#
# 1. On the "peds" sheet, stamp a new "Test" value (shared string) down
#    column AY for every data row (2-18) under the existing "Schedule"
#    header in AY1, and move that sheet's selection to AY24.
# 2. Make "adults" the active/selected sheet and set its selection to
#    AX2:AX12.

$wb = $excel.ActiveWorkbook

$peds = $wb.Worksheets.Item("peds")
$adults = $wb.Worksheets.Item("adults")

# --- peds: fill column AY ("Test") for data rows 2-18 ---
for ($r = 2; $r -le 18; $r++) {
    $peds.Range("AY$r").Value = "Test"
}

# peds keeps the selection parked on AY24 (no longer the active tab)
$peds.Range("AY24").Select()

# --- adults becomes the active sheet, with its own selection ---
$adults.Activate()
$adults.Range("AX2:AX12").Select()
